# Refresh the crypto price/volume snapshot table (columns D = Price,
# E = Volume(1h)) with the latest scraped values.
# Mirrors commit: "Updated cryptos list ... with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="26.959.81"; E="  +0.28%  "},
    @{Row=3; D="1.554.28"; E="  +0.52%  "},
    @{Row=4; E="  -0.17%  "},
    @{Row=5; E="  +0.66%  "},
    @{Row=6; E="  +0.36%  "},
    @{Row=7; E="  -0.20%  "},
    @{Row=8; D="22.13"; DText=$true; E="  +3.98%  "},
    @{Row=9; E="  -0.04%  "},
    @{Row=10; E="  +0.91%  "},
    @{Row=11; E="  +0.05%  "},
    @{Row=12; D="1.776.07"; E="  +0.53%  "},
    @{Row=13; D="1.555.38"; E="  +0.43%  "},
    @{Row=14; E="  +1.44%  "},
    @{Row=15; E="  +1.54%  "},
    @{Row=16; D="26.967.75"},
    @{Row=17; D="61.69"; DText=$true; E="  +0.24%  "},
    @{Row=18; D="217.90"; DText=$true; E="  +2.04%  "},
    @{Row=19; E="  +2.42%  "},
    @{Row=20; D="7.29"; DText=$true; E="  +1.54%  "},
    @{Row=21; E="  -0.16%  "},
    @{Row=22; E="  +1.29%  "},
    @{Row=23; E="  +0.58%  "},
    @{Row=24; E="  +0.57%  "},
    @{Row=25; D="154.46"; DText=$true; E="  +1.02%  "},
    @{Row=26; D="6.63"; DText=$true; E="  -0.13%  "},
    @{Row=27; E="  +0.83%  "},
    @{Row=28; E="  +1.17%  "},
    @{Row=29; E="  -0.18%  "},
    @{Row=30; E="  +2.15%  "},
    @{Row=31; E="  -0.53%  "},
    @{Row=32; E="  +0.56%  "},
    @{Row=33; D="1.425.95"; E="  +4.64%  "},
    @{Row=34; E="  +4.60%  "},
    @{Row=35; D="1.58"; DText=$true; E="  +3.44%  "},
    @{Row=36; D="0.978"; DText=$true; E="  +1.55%  "},
    @{Row=37; E="  +0.21%  "},
    @{Row=38; E="  +0.66%  "},
    @{Row=39; D="0.522"; DText=$true; E="  +0.80%  "},
    @{Row=40; E="  +0.79%  "},
    @{Row=41; D="5.77"; DText=$true; E="  +3.41%  "},
    @{Row=42; E="  -0.16%  "},
    @{Row=43; E="  +4.64%  "},
    @{Row=44; E="  +0.33%  "},
    @{Row=45; D="64.30"; DText=$true; E="  +1.44%  "},
    @{Row=46; E="  +1.70%  "},
    @{Row=47; D="1.690.12"; E="  +0.50%  "},
    @{Row=48; D="87.71"; DText=$true},
    @{Row=49; E="  +2.86%  "},
    @{Row=50; E="  +3.62%  "},
    @{Row=51; D="0.0955"; DText=$true; E="  +0.72%  "}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($u.ContainsKey("DText")) {
            # Force plain-text storage so values such as "22.13" or
            # "7.29" are not reinterpreted as numbers (matches the
            # original file, where every Price cell is an inline string).
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
